$wb = $excel.ActiveWorkbook

# --- 1. Add the new "en corrección" status to the "estado" lookup table (sheet "tablas") ---
$wsTablas = $wb.Worksheets.Item("tablas")
$loEstado = $wsTablas.ListObjects.Item("estado")
$newRow = $loEstado.ListRows.Add()
$newRow.Range.Item(1).Value = "en corrección"

# --- 2. Apply the new status to the two rows that are now "en corrección" ---
$wsTemplates = $wb.Worksheets.Item("templates")
$wsTemplates.Range("E3").Value = "en corrección"
$wsTemplates.Range("G3").Value = "en corrección"
$wsTemplates.Range("E4").Value = "en corrección"
$wsTemplates.Range("G4").Value = "en corrección"

# --- 3. Add a conditional formatting rule that highlights "en corrección" cells ---
$rng = $wsTemplates.Range("A1:XFD1048576")
$fc = $rng.FormatConditions.Add(9, 0, $null, $null, "en corrección")
$fc.Font.Color = 7039480
$fc.Interior.Color = 12709465
$fc.SetFirstPriority()

# --- 4. Column width tweaks on "templates" ---
$wsTemplates.Columns.Item(5).ColumnWidth = 11.28515625
$wsTemplates.Columns.Item(7).ColumnWidth = 11.28515625
$wsTemplates.Columns.Item(10).ColumnWidth = 12.140625

# --- 5. Selection / view state ---
$wsTemplates.Range("H4").Select()
$wsTablas.Range("C9").Select()
